# Automatische test-sync: 2025-08-01 23:47:50
#
# Adds a new incoming-mail row (#6, about EcoPro-700 stock / "Productinformatie"
# category) to the "Logs" sheet, bumps the matching "Productinformatie" tally
# on the "Dashboard" sheet, and extends the conditional formatting ranges and
# the dashboard bar chart's category/value series so they include the new row.

$wb  = $excel.ActiveWorkbook
$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs!A11:J11 -----------------------------------------------------
$logs.Cells.Item(11, 1).Value  = "Hebben we EcoPro-700 nog op voorraad?"
$logs.Cells.Item(11, 2).Value  = "mailmind.test@zohomail.eu"
$logs.Cells.Item(11, 3).Value  = "Testmail #6: Hebben we EcoPro-700 nog op voorraad?"
$logs.Cells.Item(11, 4).Value  = "Productinformatie"
$logs.Cells.Item(11, 5).Value  = "Geachte klant,`nHartelijk dank voor uw e-mail. Op dit moment hebben we nog EcoPro-700 op voorraad. U kunt deze bestellen via onze website of neem contact met ons op voor verdere assistentie.`nMet vriendelijke groet,`n[Bedrijfsnaam]"
$logs.Cells.Item(11, 6).Value  = "2025-08-01 23:47:25"
$logs.Cells.Item(11, 7).Value  = "Ja"
$logs.Cells.Item(11, 8).Value  = "Nee"
$logs.Cells.Item(11, 9).Value  = "Ja"
$logs.Cells.Item(11, 10).Value = "Nee"

# --- extend conditional formatting ranges on Logs from row 10 to row 11
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D11"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G11"))
$logs.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H11"))
$logs.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I11"))
$logs.Range("J2:J10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J11"))

# --- Dashboard!A5:B5 ----------------------------------------------------
$dash.Cells.Item(5, 1).Value = "Productinformatie"
$dash.Cells.Item(5, 2).Value = 1

# --- extend the bar chart's category/value series ranges to row 5 -----
$chartSeries = $dash.ChartObjects().Item(1).Chart.SeriesCollection(1)
$chartSeries.XValues = "='Dashboard'!`$A`$2:`$A`$5"
$chartSeries.Values  = "='Dashboard'!`$B`$2:`$B`$5"
